$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 77.35714
$ws.Range("I9").Value = 82
$ws.Range("J9").Value = 60.333332
$ws.Range("K9").Value = 82
$ws.Range("L9").Value = 60.333332
$ws.Range("M9").Value = 87
$ws.Range("N9").Value = -398.333332

$ws.Range("H100").Value = 7659.6665
$ws.Range("I100").Value = 2626.625
$ws.Range("J100").Value = 14980.454
$ws.Range("K100").Value = 2626.625
$ws.Range("L100").Value = 14980.454
$ws.Range("M100").Value = -2085.625
$ws.Range("N100").Value = -16062.454

$ws.Range("H116").Value = 15101.833
$ws.Range("I116").Value = 7123
$ws.Range("J116").Value = 20801
$ws.Range("K116").Value = 7123
$ws.Range("L116").Value = 20801
$ws.Range("M116").Value = -3681
$ws.Range("N116").Value = -27685

$ws.Range("H135").Value = 3551.4443
$ws.Range("I135").Value = 3352.3572
$ws.Range("J135").Value = 4248.25
$ws.Range("K135").Value = 30171.2148
$ws.Range("L135").Value = 38234.25
$ws.Range("M135").Value = -27636.2148
$ws.Range("N135").Value = -43304.25

$ws.Range("H138").Value = 5413.726
$ws.Range("J138").Value = 6804.9473
$ws.Range("L138").Value = 20414.8419
$ws.Range("N138").Value = -30694.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3814.1287
$ws.Range("I32").Value = 3283.4478
$ws.Range("J32").Value = 15666
$ws.Range("K32").Value = 3283.4478
$ws.Range("L32").Value = 15666
$ws.Range("M32").Value = -2996.4478
$ws.Range("N32").Value = -16240

$ws.Range("H61").Value = 4334.55
$ws.Range("I61").Value = 4331.871
$ws.Range("K61").Value = 4331.871
$ws.Range("M61").Value = -4119.871

$ws.Range("H102").Value = 1249.3422
$ws.Range("I102").Value = 1179.8889
$ws.Range("K102").Value = 1179.8889
$ws.Range("M102").Value = 442.1111000000001

$ws.Range("H132").Value = 4053.9207
$ws.Range("I132").Value = 3350.2144
$ws.Range("J132").Value = 5461.3335
$ws.Range("K132").Value = 10050.6432
$ws.Range("L132").Value = 16384.0005
$ws.Range("M132").Value = -7520.643199999999
$ws.Range("N132").Value = -21444.0005

$ws.Range("H136").Value = 4334.55
$ws.Range("I136").Value = 4331.871
$ws.Range("K136").Value = 12995.613
$ws.Range("M136").Value = -10445.613

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J20").Value = 3068.6365
$ws.Range("L20").Value = 3068.6365
$ws.Range("N20").Value = -3562.6365

$ws.Range("H22").Value = 246.27272
$ws.Range("I22").Value = 171.1
$ws.Range("K22").Value = 171.1
$ws.Range("M22").Value = 1.900000000000006

$ws.Range("H99").Value = 3678.7
$ws.Range("I99").Value = 3590.2222
$ws.Range("K99").Value = 3590.2222
$ws.Range("M99").Value = -2092.2222

$ws.Range("H134").Value = 31996.324
$ws.Range("I134").Value = 4914.56
$ws.Range("J134").Value = 88416.664
$ws.Range("K134").Value = 14743.68
$ws.Range("L134").Value = 265249.992
$ws.Range("M134").Value = -12208.68
$ws.Range("N134").Value = -270319.992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 64592
$ws.Range("J97").Value = 64592
$ws.Range("L97").Value = 64592
$ws.Range("N97").Value = -66574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1022.55554
$ws.Range("I18").Value = 673.36365
$ws.Range("J18").Value = 2559
$ws.Range("K18").Value = 2020.09095
$ws.Range("L18").Value = 7677
$ws.Range("M18").Value = -1851.09095
$ws.Range("N18").Value = -8015

$ws.Range("H121").Value = 28611934
$ws.Range("I121").Value = 809
$ws.Range("K121").Value = 2427
$ws.Range("M121").Value = -1117

$ws.Range("H132").Value = 396786.97
$ws.Range("J132").Value = 593153.9399999999
$ws.Range("L132").Value = 5338385.459999999
$ws.Range("N132").Value = -5343445.459999999

$ws.Range("H139").Value = 5047.5527
$ws.Range("I139").Value = 1445.4
$ws.Range("K139").Value = 4336.200000000001
$ws.Range("M139").Value = 803.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 30547.795
$ws.Range("I132").Value = 6169.2144
$ws.Range("K132").Value = 18507.6432
$ws.Range("M132").Value = -15977.6432

$ws.Range("H135").Value = 90977270
$ws.Range("J135").Value = 90977270
$ws.Range("L135").Value = 90977270
$ws.Range("N135").Value = -90987410

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 678.5
$ws.Range("I16").Value = 678.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 678.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -508.5
$ws.Range("N16").ClearContents()

$ws.Range("H61").Value = 3946
$ws.Range("I61").Value = 4596.9165
$ws.Range("J61").Value = 1993.25
$ws.Range("K61").Value = 4596.9165
$ws.Range("L61").Value = 1993.25
$ws.Range("M61").Value = -4394.9165
$ws.Range("N61").Value = -2397.25

$ws.Range("H68").Value = 45192.043
$ws.Range("I68").Value = 1676.8235
$ws.Range("J68").Value = 168485.17
$ws.Range("K68").Value = 1676.8235
$ws.Range("L68").Value = 168485.17
$ws.Range("M68").Value = -927.8235
$ws.Range("N68").Value = -169983.17

$ws.Range("H71").Value = 45192.043
$ws.Range("I71").Value = 1676.8235
$ws.Range("J71").Value = 168485.17
$ws.Range("K71").Value = 8384.1175
$ws.Range("L71").Value = 842425.8500000001
$ws.Range("M71").Value = -4640.1175
$ws.Range("N71").Value = -849913.8500000001

$ws.Range("H113").Value = 3946
$ws.Range("I113").Value = 4596.9165
$ws.Range("J113").Value = 1993.25
$ws.Range("K113").Value = 4596.9165
$ws.Range("L113").Value = 1993.25
$ws.Range("M113").Value = -2426.9165
$ws.Range("N113").Value = -6333.25

$ws.Range("H136").Value = 150634.55
$ws.Range("I136").Value = 236585.28
$ws.Range("J136").Value = 8485.27
$ws.Range("K136").Value = 709755.84
$ws.Range("L136").Value = 25455.81
$ws.Range("M136").Value = -707205.84
$ws.Range("N136").Value = -30555.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6599.0415
$ws.Range("I81").Value = 2771.6155
$ws.Range("J81").Value = 11122.363
$ws.Range("K81").Value = 5543.231
$ws.Range("L81").Value = 22244.726
$ws.Range("M81").Value = -4482.231
$ws.Range("N81").Value = -24366.726

$ws.Range("H84").Value = 6599.0415
$ws.Range("I84").Value = 2771.6155
$ws.Range("J84").Value = 11122.363
$ws.Range("K84").Value = 27716.155
$ws.Range("L84").Value = 111223.63
$ws.Range("M84").Value = -22412.155
$ws.Range("N84").Value = -121831.63

$ws.Range("H132").Value = 17645.713
$ws.Range("I132").Value = 2082.2222
$ws.Range("K132").Value = 6246.6666
$ws.Range("M132").Value = -3716.6666

$ws.Range("H136").Value = 275584.3
$ws.Range("I136").Value = 306011.1
$ws.Range("J136").Value = 184304
$ws.Range("K136").Value = 918033.2999999999
$ws.Range("L136").Value = 552912
$ws.Range("M136").Value = -915483.2999999999
$ws.Range("N136").Value = -558012
